# Apply the table style change on the table on slide 16
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table
$tbl.ApplyStyle("{C70F4F56-C206-46C3-8ACA-FAF2FFE29510}")

# Swap the presentation's active theme color scheme from "Integral" colours
# to the "Office Theme" colours (the diff swaps the contents of theme1.xml
# and theme2.xml; theme2.xml is the part actually wired up to the slide
# master/presentation, so re-pointing its colour scheme reproduces the
# visible effect of the swap).
$tcs = $s.ThemeColorScheme
$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
